# heaveDecayCases.xlsx - "updated heave decay but results still look bad"
#
# 1. Header cell A2 ("m") is renamed to "heaveDecay".
# 2. Column A is widened to fit the new header text (customWidth).
# 3. The shared-formula block in B3:G7 (=-0.04394+$A$n) is replaced by
#    plain literal values equal to column A's value for that row - i.e.
#    the "decay" columns no longer compute an offset, they just mirror
#    the initial displacement in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Rename the header in A2.
$ws.Range("A2").Value = "heaveDecay"

# 2. Resize column A so the new, longer header text fits.
$ws.Columns.Item(1).ColumnWidth = 10.75

# 3. Replace the formulas in B3:G7 with the literal value taken from
#    column A of the same row (the formula result is no longer derived).
foreach ($row in 3..7) {
    $aValue = $ws.Cells.Item($row, 1).Value2
    $ws.Range("B$row`:G$row").Value = $aValue
}
